$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New development task rows (Ing / Tiempo / Actividad / Descripcion)
$ws.Range("B6").Value = "?"
$ws.Range("C6").Value = "1 dias"
$ws.Range("E6").Value = "Sistema de inicio de sesion"
$ws.Range("G6").Value = "Los datos seran administrados y solo autorizados por los duenos o creadores de la cuenta"

$ws.Range("B7").Value = "?"
$ws.Range("C7").Value = "5 dias"
$ws.Range("E7").Value = "Dashboard"
$ws.Range("G7").Value = "Toda la informacion del negocio, desde datos, compras, ventas, ganancias, perdidas"

$ws.Range("B8").Value = "?"
$ws.Range("C8").Value = "3 dias"
$ws.Range("E8").Value = "Inventario"
$ws.Range("G8").Value = "Apartado donde el usuario podra ver, revisar, administrar sus productos, dar de alta, y revisar inventario"

$ws.Range("B9").Value = "?"
$ws.Range("C9").Value = "5 dias"
$ws.Range("E9").Value = "Shop"
$ws.Range("G9").Value = "Tienda donde se cobraran los productos y se venderan al cliente"

# C8/C9 were previously blank cells without an explicit style; copy the
# existing column C formatting (center/middle aligned, wrapped text) onto
# them so they match the rest of the "Tiempo" column.
$ws.Range("C3").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection as it was left after the edit
$ws.Range("G13").Select()
